$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# The paragraph "* What is diference between React and ReactDOM" contains a
# misspelling ("diference") that Word had flagged with proofing marks
# (<w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>).
# The fix corrects the spelling to "difference". Because the word is now
# spelled correctly, the spell-check markers around it are removed, and the
# correction (typing the missing "f") leaves the surrounding text split
# into three runs: "... What is dif" / "f" / "erence between ...".
# --------------------------------------------------------------------------

# Step 1: correct the typo. Replacing the whole "What is diference between"
# phrase (not just the misspelled word) means the replace operation spans
# across the <w:proofErr> boundaries on both sides of "diference", so Word
# drops those now-unneeded proofing marks at the same time as fixing the text.
$d.Content.Find.Execute("What is diference between", $false, $false, $false, $false, $false, $true, 1, $false, "What is difference between", 2)

# Step 2: locate the corrected word so we can split out the single
# character that was "typed in" to fix the typo (the second "f" of
# "difference") into its own run, matching how Word records an in-place
# correction.
$rng = $d.Content
$rng.Find.Text = "difference between React and"
$rng.Find.Execute()

$wordStart = $rng.Start          # start of "difference"
$fStart = $wordStart + 3         # "d","i","f","f" -> index 3 is the 2nd "f"
$fEnd = $fStart + 1

$fixedChar = $d.Range($fStart, $fEnd)

# Toggling a character-formatting property and then reverting it forces the
# run containing this single character to be split off from its neighbours
# without altering the visible formatting of the text.
$fixedChar.Bold = 1
$fixedChar.Bold = 0
